# Insert a new data row at row 61 (pushing existing rows 61-134 down to 62-135)
# and populate it with the new record described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 61, shifting rows 61..134 down to 62..135.
$ws.Rows.Item(61).Insert()

# Fill in the new row 61 with the new record's data.
$ws.Range("A61").Value = 1
$ws.Range("B61").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C61").Value = "Arica y Parinacota"
$ws.Range("D61").Value = (Get-Date -Year 2022 -Month 10 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E61").Value = 15
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100106
$ws.Range("H61").Value = "Oleaginosos"
$ws.Range("I61").Value = 100106002
$ws.Range("J61").Value = "Palta"
$ws.Range("K61").Value = "Hass"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 600
$ws.Range("N61").Value = 22000
$ws.Range("O61").Value = 23000
$ws.Range("P61").Value = 22500
$ws.Range("Q61").Value = "`$/bandeja 10 kilos"
$ws.Range("R61").Value = "Perú"
$ws.Range("S61").Value = 2250
$ws.Range("T61").Value = 10
